# Update login/sign-up/profile-pic schema on the DATA sheet, and refresh
# the corresponding selection/zoom UI state.

$wb = $excel.ActiveWorkbook
$data = $wb.Worksheets.Item("DATA")

# --- Header row (row 1): columns that keep re-using already-known text ---
$data.Cells.Item(1, 1).Value = "user_id"
$data.Cells.Item(1, 2).Value = "username"
$data.Cells.Item(1, 3).Value = "role"
$data.Cells.Item(1, 4).Value = "name"
$data.Cells.Item(1, 5).Value = "password"
$data.Cells.Item(1, 7).Value = "last_login"

# --- Data row (row 2): user_id/username/password are unchanged values; ---
# --- role/name are rewritten with their new display text (new strings). ---
$data.Cells.Item(2, 1).Value = 1
$data.Cells.Item(2, 2).Value = "napatswift"
$data.Cells.Item(2, 3).Value = "ADMIN"
$data.Cells.Item(2, 4).Value = "Napat"
$data.Cells.Item(2, 5).Value = "NAPAT1"
$data.Cells.Item(2, 6).ClearContents()

# --- New header cells for the added columns ---
$data.Cells.Item(1, 6).Value = "picturePath"
$data.Cells.Item(1, 8).Value = "isBanned"
$data.Cells.Item(1, 9).Value = "loginAttempt"
$data.Cells.Item(1, 10).Value = "hasStore"
$data.Cells.Item(1, 11).Value = "store"

# --- New data cells for the added columns ---
$data.Cells.Item(2, 8).Value = $false
$data.Cells.Item(2, 10).Value = $false

# --- View state: zoom + selection for the DATA sheet ---
$data.Activate()
$excel.ActiveWindow.Zoom = 157
$data.Range("H4").Select()
